$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Insert a new bold "Play Cash Eruption slot for free" paragraph
#    right before the final paragraph (the one that currently reads
#    "Create an attention-grabbing feature image for ...").
# -----------------------------------------------------------------
$lastIdx  = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIdx)

$lastPara.Range.InsertParagraphBefore()

$newIdx  = $lastIdx
$newPara = $d.Paragraphs($newIdx)

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cash Eruption slot for free</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($xmlFrag)

# -----------------------------------------------------------------
# 2. Swap the old "Create an attention-grabbing feature image ..."
#    italic paragraph text for the meta-description sentence (minus
#    the "Meta description: " prefix) while keeping its italics.
# -----------------------------------------------------------------
$finalIdx  = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($finalIdx)
$finalPara.Range.Find.Execute(
    "Create an attention-grabbing feature image for the Cash Eruption slot game in a cartoon style. The image should feature a happy Maya warrior with glasses. The warrior can be holding a bag of gold coins and standing next to a volcano that is erupting fire and cash. Use bright colors and bold outlines to make the image pop and convey the excitement of the game. Make sure the text " + [char]34 + "Cash Eruption" + [char]34 + " is prominently displayed in the image along with any other relevant information, such as the game developer or where it can be played.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Cash Eruption, the Aztec-themed slot machine game with four jackpots and bonus rounds available to play for free.",
    2)

# -----------------------------------------------------------------
# 3. Remove the original "Meta description: Read our review ..."
#    paragraph that used to sit right under the page's H1 heading.
# -----------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()
